$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Insert a new "properties" worksheet between "survey" and "settings"
# ------------------------------------------------------------------
$surveySheet = $wb.Worksheets.Item("survey")
$propsSheet = $wb.Worksheets.Add($null, $surveySheet)
$propsSheet.Name = "properties"

# ------------------------------------------------------------------
# Populate the properties sheet. Values are written in the specific
# order below so the resulting shared-string table ordering matches
# the original authoring session.
# ------------------------------------------------------------------
$propsSheet.Range("A1").Value = "partition"
$propsSheet.Range("B1").Value = "aspect"
$propsSheet.Range("C1").Value = "key"
$propsSheet.Range("D1").Value = "type"
$propsSheet.Range("E1").Value = "value"

$propsSheet.Range("A2").Value = "Table"
$propsSheet.Range("B2").Value = "default"
$propsSheet.Range("C2").Value = "colOrder"
$propsSheet.Range("D2").Value = "array"
$propsSheet.Range("E2").Value = '["Description","Image_uriFragment","Location_latitude","Location_longitude","Date_and_Time"]'

$propsSheet.Range("A3").Value = "Table"
$propsSheet.Range("B3").Value = "default"
$propsSheet.Range("C3").Value = "defaultViewType"
$propsSheet.Range("D3").Value = "string"
$propsSheet.Range("E3").Value = "MAP"

$propsSheet.Range("A4").Value = "Table"
$propsSheet.Range("B4").Value = "default"
$propsSheet.Range("C4").Value = "detailViewFileName"
$propsSheet.Range("E4").Value = "config/tables/geotagger/html/geo_detail.html"

$propsSheet.Range("C7").Value = "keyColorRuleType"
$propsSheet.Range("D7").Value = "string"
$propsSheet.Range("E7").Value = "None"

$propsSheet.Range("A5").Value = "Table"
$propsSheet.Range("B5").Value = "default"
$propsSheet.Range("C5").Value = "listViewFileName"

$propsSheet.Range("D4").Value = "configpath"
$propsSheet.Range("D5").Value = "configpath"
$propsSheet.Range("E5").Value = "config/tables/geotagger/html/geo_list.html"

$propsSheet.Range("A6").Value = "Table"
$propsSheet.Range("B6").Value = "default"
$propsSheet.Range("C6").Value = "mapListViewFileName"

$propsSheet.Range("A7").Value = "TableMapFragment"
$propsSheet.Range("B7").Value = "default"

$propsSheet.Range("D6").Value = "configpath"
$propsSheet.Range("E6").Value = "config/tables/geotagger/html/geo_list.html"

$propsSheet.Range("A8").Value = "TableMapFragment"
$propsSheet.Range("B8").Value = "default"
$propsSheet.Range("C8").Value = "keyMapLatCol"
$propsSheet.Range("D8").Value = "string"
$propsSheet.Range("E8").Value = "Location_latitude"

$propsSheet.Range("A9").Value = "TableMapFragment"
$propsSheet.Range("B9").Value = "default"
$propsSheet.Range("C9").Value = "keyMapLongCol"
$propsSheet.Range("D9").Value = "string"
$propsSheet.Range("E9").Value = "Location_longitude"

# ------------------------------------------------------------------
# Column widths for the new sheet
# ------------------------------------------------------------------
$propsSheet.Columns.Item(1).ColumnWidth = 19.944010416666668
$propsSheet.Columns.Item(2).ColumnWidth = 13.166666666666666
$propsSheet.Columns.Item(3).ColumnWidth = 18.944010416666668
$propsSheet.Columns.Item(4).ColumnWidth = 13.166666666666666
$propsSheet.Columns.Item(5).ColumnWidth = 13.166666666666666

# ------------------------------------------------------------------
# Selection on the new properties sheet (becomes the active sheet)
# ------------------------------------------------------------------
[void]$propsSheet.Range("B16").Select()
